$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 0
$ws1.Range("F4").Value = 19527
$ws1.Range("F5").Value = 782
$ws1.Range("F7").Value = 0
$ws1.Range("F9").Value = 7392
$ws1.Range("F10").Value = 483
$ws1.Range("F12").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 101
$ws1.Range("F16").Value = 0
$ws1.Range("F19").Value = 363
$ws1.Range("F22").Value = 45
$ws1.Range("F23").Value = 0
$ws1.Range("F26").Value = 1068
$ws1.Range("F27").Value = 0
$ws1.Range("F28").Value = 0
$ws1.Range("F30").Value = 5219
$ws1.Range("F32").Value = 0
$ws1.Range("F33").Value = 2568
$ws1.Range("F34").Value = 0
$ws1.Range("F37").Value = 12489
$ws1.Range("F38").Value = 1318
$ws1.Range("F40").Value = 13
$ws1.Range("F41").Value = 0
$ws1.Range("F42").Value = 248
$ws1.Range("F43").Value = 0
$ws1.Range("F45").Value = 0
$ws1.Range("F46").Value = 0

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 19527
$ws4.Range("F5").Value = 782
$ws4.Range("F7").Value = 0
$ws4.Range("F9").Value = 7392
$ws4.Range("F10").Value = 483
$ws4.Range("F12").Value = 248
$ws4.Range("F13").Value = 0
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 45
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 0
$ws4.Range("F25").Value = 307
$ws4.Range("F26").Value = 1068
$ws4.Range("F28").Value = 8
$ws4.Range("F29").Value = 0
$ws4.Range("F30").Value = 5219
$ws4.Range("F31").Value = 0
$ws4.Range("F33").Value = 0
$ws4.Range("F35").Value = 2571
$ws4.Range("F38").Value = 0
$ws4.Range("F39").Value = 12489
$ws4.Range("F40").Value = 0
$ws4.Range("F42").Value = 0
$ws4.Range("F44").Value = 248
$ws4.Range("F45").Value = 335
$ws4.Range("F46").Value = 0
$ws4.Range("F48").Value = 0

Write-Output "Applied all F-column updates"
